# UPDATE LAYOUT MASTER BARANG SATU SATUAN
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove now-unused Sheet2 / Sheet3 ---
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

# --- Drop the trailing, unused K:L columns on Sheet1 ---
$ws.Range("K1:L5").EntireColumn.Delete()

# --- Rebuild header row (A1:J1) ---
$ws.Range("A1").Value = "NAMABARANG"
$ws.Range("B1").Value = "KATEGORI"
$ws.Range("C1").Value = "SUBKATEGORI"
$ws.Range("D1").Value = "SATUAN"
$ws.Range("E1").Value = "HARGABELI"
$ws.Range("F1").Value = "HARGAJUAL"
$ws.Range("G1").Value = "JMLBARANG"
$ws.Range("H1").Value = "BARCODE"
$ws.Range("I1").Value = "NAMASUPPLIER"
$ws.Range("J1").Value = "RAK"

# --- Populate data rows (A2:J5) ---
$ws.Range("A2").Value = "KOPI"
$ws.Range("B2").Value = "MINUMAN"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "PCS"
$ws.Range("E2").Value = 5000
$ws.Range("F2").Value = 6500
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = "32140P"
$ws.Range("I2").Value = "UDI JAYA"
$ws.Range("J2").Value = "RAK 1"

$ws.Range("A3").Value = "TEH JAWA"
$ws.Range("B3").Value = "MINUMAN"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "PCS"
$ws.Range("E3").Value = 7500
$ws.Range("F3").Value = 9000
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = "H09327392"
$ws.Range("I3").Value = "BUDI LUHUR"
$ws.Range("J3").Value = "RAK 2"

$ws.Range("A4").Value = "GULAKU"
$ws.Range("B4").Value = "BAHAN MASAKAN"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "KG"
$ws.Range("E4").Value = 15000
$ws.Range("F4").Value = 17000
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = "P4387439"
$ws.Range("J4").Value = "RAK 3"

$ws.Range("A5").Value = "MAKARONI"
$ws.Range("B5").Value = "BAHAN MASAKAN"
$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = "KG"
$ws.Range("E5").Value = 25000
$ws.Range("F5").Value = 28000
$ws.Range("G5").Value = 50
$ws.Range("H5").Value = "J0923742"
$ws.Range("J5").Value = "-"

# --- Selection / view state matches the committed sheet ---
$ws.Range("J6").Select()
